$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark from its old location (end
#    of the "Começo do uso do DTO..." paragraph). It gets re-created
#    later, inside the new "Uso da paginação..." paragraph.
# ---------------------------------------------------------------------
try {
    $oldGoBack = $d.Bookmarks("_GoBack")
    $oldGoBack.Delete()
} catch {
}

# ---------------------------------------------------------------------
# 2. Materialize the built-in "Hyperlink" character style cleanly (no
#    w:customStyle, correct w:themeColor) by applying it by name to a
#    throw-away run, tweak its properties to match the target, then
#    remove the throw-away paragraph again. The style definition stays
#    behind in styles.xml once it has been used at least once.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$tempPara = $d.Paragraphs.Last
$tempRange = $tempPara.Range
$tempRange.Collapse(0)
$tempRange.InsertAfter("x")
$tempRunRange = $d.Range($tempRange.Start, $tempRange.Start + 1)
$tempRunRange.Style = "Hyperlink"

$hlStyle = $d.Styles("Hyperlink")
$hlStyle.BaseStyle = $d.Styles("Fontepargpadro")
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true

$tempPara.Range.Delete()

# ---------------------------------------------------------------------
# 3. Append a new, empty paragraph after the last paragraph of the
#    document; this becomes our insertion host so the existing last
#    paragraph's own content ("... quer.") is left completely intact.
# ---------------------------------------------------------------------
$lastPara2 = $d.Paragraphs.Last
$lastPara2.Range.InsertParagraphAfter()
$hostPara = $d.Paragraphs.Last
$hostRange = $hostPara.Range

# ---------------------------------------------------------------------
# 4. Insert the six new paragraphs as raw OOXML. The final paragraph is
#    left empty on purpose: InsertXML leaves the host paragraph mark
#    behind as a trailing empty paragraph, matching the diff's final
#    blank paragraph.
# ---------------------------------------------------------------------
$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pBdr>
      <w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/>
    </w:pBdr>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve">No caso na hora de mostrar o </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t>endpoint</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve"> Categorias, n&#227;o abrir&#225; o &#8220;leque&#8221; dos produtos</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve">Colocando </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t>endpoint</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve"> de pagina&#231;&#227;o</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t>Uso da pagina&#231;&#227;o</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve"> com par&#226;metros na requisi&#231;&#227;o</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve">: busca de tanto em tanto no banco de dados, para n&#227;o sobrecarregar o sistema. </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t>Ex</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t>: buscar de 20 em 20.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve">Quando criamos o m&#233;todo do </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t>page</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve"> n&#227;o vamos fazer assim: /categorias/</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t>page</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve">/0/20 (p&#225;gina 0 com 20 linhas), ou seja, n&#227;o ser&#225; </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t>vari&#225;veis</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve"> do pr&#243;prio path, e sim como par&#226;metros...</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
      <w:color w:val="505050"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>
    <w:t xml:space="preserve">Teste da pagina&#231;&#227;o: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
      <w:color w:val="505050"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>http://localhost:8080/categorias/page?linesPerPage=3&amp;page=1</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$hostRange.InsertXML($xmlFrag)

# ---------------------------------------------------------------------
# 5. Re-create the "_GoBack" bookmark right after "...na requisição"
#    and before ": busca de tanto em tanto...", matching its new
#    position in the target document.
# ---------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("com parâmetros na requisição", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bmPoint = $d.Range($findRange.End, $findRange.End)
    $d.Bookmarks.Add("_GoBack", $bmPoint)
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
